# Update cryptos list: prices and 1h volume percentages refreshed,
# plus a few coins swapped ranking position (rows 22/23, 42/43, 45/46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks numeric need to be forced to
# Text format first, otherwise Excel would auto-convert them to numbers
# (the source data keeps them as text, e.g. "579.37", "0.999", "20.49").
$textFormatCells = @("D5", "D6", "D7", "D9", "D11", "D14", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D30", "D31", "D32", "D33", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D49", "D50")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "61.035.42"
$ws.Range("E2").Value = "  +0.69%  "
# Row 3
$ws.Range("D3").Value = "2.673.82"
$ws.Range("E3").Value = "  +2.55%  "
# Row 4
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("D5").Value = "579.37"
$ws.Range("E5").Value = "  +0.15%  "
# Row 6
$ws.Range("D6").Value = "145.17"
$ws.Range("E6").Value = "  +1.43%  "
# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.21%  "
# Row 8
$ws.Range("E8").Value = "  -0.04%  "
# Row 9
$ws.Range("D9").Value = "6.61"
$ws.Range("E9").Value = "  +1.48%  "
# Row 10
$ws.Range("E10").Value = "  +1.46%  "
# Row 11
$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  +5.58%  "
# Row 12
$ws.Range("E12").Value = "  +1.17%  "
# Row 13
$ws.Range("D13").Value = "3.144.48"
$ws.Range("E13").Value = "  +2.18%  "
# Row 14
$ws.Range("D14").Value = "25.93"
$ws.Range("E14").Value = "  +11.51%  "
# Row 15
$ws.Range("D15").Value = "61.012.22"
$ws.Range("E15").Value = "  +0.54%  "
# Row 16
$ws.Range("E16").Value = "  +1.39%  "
# Row 17
$ws.Range("D17").Value = "2.664.38"
$ws.Range("E17").Value = "  +1.57%  "
# Row 18
$ws.Range("D18").Value = "11.68"
$ws.Range("E18").Value = "  +3.33%  "
# Row 19
$ws.Range("E19").Value = "  +2.50%  "
# Row 20
$ws.Range("D20").Value = "351.23"
$ws.Range("E20").Value = "  +0.77%  "
# Row 21
$ws.Range("D21").Value = "6.94"
$ws.Range("E21").Value = "  +0.37%  "
# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.31%  "
# Row 23
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "0.537"
$ws.Range("E23").Value = "  +2.79%  "
# Row 24
$ws.Range("E24").Value = "  +0.99%  "
# Row 25
$ws.Range("D25").Value = "0.162"
$ws.Range("E25").Value = "  +1.03%  "
# Row 26
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.27%  "
# Row 27
$ws.Range("E27").Value = "  +5.26%  "
# Row 28
$ws.Range("E28").Value = "  +5.76%  "
# Row 29
$ws.Range("D29").Value = "0.0₃0816"
$ws.Range("E29").Value = "  +2.99%  "
# Row 30
$ws.Range("D30").Value = "6.89"
$ws.Range("E30").Value = "  +8.05%  "
# Row 31
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.13%  "
# Row 32
$ws.Range("D32").Value = "165.50"
$ws.Range("E32").Value = "  +2.21%  "
# Row 33
$ws.Range("D33").Value = "19.98"
$ws.Range("E33").Value = "  +2.15%  "
# Row 34
$ws.Range("E34").Value = "  +11.29%  "
# Row 35
$ws.Range("E35").Value = "  +6.08%  "
# Row 36
$ws.Range("E36").Value = "  +6.06%  "
# Row 37
$ws.Range("E37").Value = "  +4.10%  "
# Row 38
$ws.Range("D38").Value = "335.42"
$ws.Range("E38").Value = "  +12.25%  "
# Row 39
$ws.Range("E39").Value = "  +4.37%  "
# Row 40
$ws.Range("D40").Value = "38.66"
$ws.Range("E40").Value = "  +1.99%  "
# Row 41
$ws.Range("D41").Value = "0.885"
$ws.Range("E41").Value = "  +4.94%  "
# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "20.49"
$ws.Range("E42").Value = "  +2.63%  "
# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "5.19"
$ws.Range("E43").Value = "  +2.99%  "
# Row 44
$ws.Range("D44").Value = "134.55"
$ws.Range("E44").Value = "  -0.77%  "
# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.100"
$ws.Range("E45").Value = "  +1.57%  "
# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0564"
$ws.Range("E46").Value = "  +3.01%  "
# Row 47
$ws.Range("E47").Value = "  +3.42%  "
# Row 48
$ws.Range("E48").Value = "  +1.71%  "
# Row 49
$ws.Range("D49").Value = "20.53"
$ws.Range("E49").Value = "  +3.62%  "
# Row 50
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.36%  "
# Row 51
$ws.Range("D51").Value = "2.099.11"
$ws.Range("E51").Value = "  +3.58%  "
